# quarterly.xlsx update: roll every quarterly table forward by one quarter.
# Columns E..N (5..14) hold ten consecutive quarters. For every table (a
# header row naming the quarters, plus its data rows) we drop the oldest
# quarter (column E), shift the remaining nine quarters left by one column,
# and populate the newly-freed last column (N) with the new quarter's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$FIRST_COL = 5   # column E
$LAST_COL  = 14  # column N

function Shift-QuarterRow($row, $newValue) {
    $vals = @()
    for ($c = $FIRST_COL; $c -le $LAST_COL; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value()
    }
    for ($c = $FIRST_COL; $c -le ($LAST_COL - 1); $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - $FIRST_COL + 1]
    }
    $ws.Cells.Item($row, $LAST_COL).Value = $newValue
}

# --- Header rows: the quarter-label row of each of the six tables -------
$headerRows = @(8, 20, 33, 46, 58, 71)
foreach ($r in $headerRows) {
    Shift-QuarterRow $r "فصل چهارم منتهی به 1401/12"
}

# --- Table 1: مقدار تولید (rows 10-16) -----------------------------------
Shift-QuarterRow 10 "-"
Shift-QuarterRow 11 0
Shift-QuarterRow 12 0
Shift-QuarterRow 13 13751
Shift-QuarterRow 14 22863
Shift-QuarterRow 15 16659
Shift-QuarterRow 16 53273

# --- Table 2: مقدار فروش (rows 22-29) -------------------------------------
Shift-QuarterRow 22 "-"
Shift-QuarterRow 23 "-"
Shift-QuarterRow 24 0
Shift-QuarterRow 25 0
Shift-QuarterRow 26 14569
Shift-QuarterRow 27 24186
Shift-QuarterRow 28 17853
Shift-QuarterRow 29 56608

# --- Table 3: مبلغ فروش (rows 35-42) --------------------------------------
Shift-QuarterRow 35 "-"
Shift-QuarterRow 36 "-"
Shift-QuarterRow 37 0
Shift-QuarterRow 38 13085
Shift-QuarterRow 39 4284127
Shift-QuarterRow 40 5330764
Shift-QuarterRow 41 2609196
Shift-QuarterRow 42 12237172

# --- Table 4: نرخ فروش (rows 48-54) ---------------------------------------
Shift-QuarterRow 48 "-"
Shift-QuarterRow 49 "-"
Shift-QuarterRow 50 "-"
Shift-QuarterRow 51 "-"
Shift-QuarterRow 52 294057725
Shift-QuarterRow 53 220407012
Shift-QuarterRow 54 146148883

# --- Table 5: مبلغ بهای تمام شده (rows 60-67) -----------------------------
Shift-QuarterRow 60 "-"
Shift-QuarterRow 61 "-"
Shift-QuarterRow 62 0
Shift-QuarterRow 63 -10376
Shift-QuarterRow 64 -2883703
Shift-QuarterRow 65 -3733122
Shift-QuarterRow 66 -2027016
Shift-QuarterRow 67 -8654217

# --- Table 6: سود ناخالص (rows 73-80) -------------------------------------
Shift-QuarterRow 73 "-"
Shift-QuarterRow 74 "-"
Shift-QuarterRow 75 0
Shift-QuarterRow 76 2709
Shift-QuarterRow 77 1400424
Shift-QuarterRow 78 1597642
Shift-QuarterRow 79 582180
Shift-QuarterRow 80 3582955
